$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edit: reservoir storage values in column B (rows 3-59) were
# re-expressed in a smaller unit (multiplied by 1e8) ---
for ($r = 3; $r -le 59; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $old = $cell.Value()
    $cell.Value = $old * 100000000
}

# --- Rows 49-59 regain the (blank) column-C cell that rows 3-48 already
# have, matching formatting, by copying the formats from C48 ---
$ws.Range("C48").Copy()
$ws.Range("C49:C59").PasteSpecial(-4122)

# --- Column C widened to fit ---
$ws.Columns.Item(3).ColumnWidth = 12

# --- Selection moved ---
$ws.Range("H30").Select()
